# Apply the latest cryptos.xlsx price/volume refresh (GitHub Actions feed update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A few "Price" values are numeric-looking text (e.g. "0.580") that Excel would
# otherwise coerce into a trimmed number on assignment; force those specific cells
# to keep their literal text formatting (trailing zero, etc.).
foreach ($addr in @("D7", "D32", "D43")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "43.002.27"
$ws.Range("E2").Value = "  -5.13%  "
$ws.Range("D3").Value = "2.219.58"
$ws.Range("E3").Value = "  -6.33%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "320.67"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").Value = "98.92"
$ws.Range("E6").Value = "  -8.67%  "
$ws.Range("D7").Value = "0.580"
$ws.Range("E7").Value = "  -8.64%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "0.556"
$ws.Range("E9").Value = "  -9.44%  "
$ws.Range("D10").Value = "36.72"
$ws.Range("E10").Value = "  -10.38%  "
$ws.Range("E11").Value = "  -3.65%  "
$ws.Range("D12").Value = "0.0825"
$ws.Range("E12").Value = "  -10.23%  "
$ws.Range("E13").Value = "  -10.25%  "
$ws.Range("E14").Value = "  -1.66%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.561.91"
$ws.Range("E15").Value = "  -6.06%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "0.862"
$ws.Range("E16").Value = "  -12.26%  "
$ws.Range("E17").Value = "  -6.95%  "
$ws.Range("D18").Value = "2.220.39"
$ws.Range("E18").Value = "  -6.27%  "
$ws.Range("D19").Value = "42.940.12"
$ws.Range("E19").Value = "  -5.11%  "
$ws.Range("D20").Value = "14.05"
$ws.Range("E20").Value = "  -9.58%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "6.54"
$ws.Range("E21").Value = "  -10.34%  "
$ws.Range("B22").Value = "ShibaInu"
$ws.Range("C22").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D22").Value = "0.0₃0962"
$ws.Range("E22").Value = "  -9.54%  "
$ws.Range("D23").Value = "64.94"
$ws.Range("E23").Value = "  -11.32%  "
$ws.Range("D24").Value = "3.18"
$ws.Range("E24").Value = "  -11.74%  "
$ws.Range("D25").Value = "236.15"
$ws.Range("E25").Value = "  -10.67%  "
$ws.Range("E26").Value = "  -7.75%  "
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("D28").Value = "4.03"
$ws.Range("E28").Value = "  +1.09%  "
$ws.Range("D29").Value = "2.24"
$ws.Range("E29").Value = "  -2.81%  "
$ws.Range("D30").Value = "9.95"
$ws.Range("E30").Value = "  -11.33%  "
$ws.Range("E31").Value = "  -16.00%  "
$ws.Range("D32").Value = "35.60"
$ws.Range("E32").Value = "  -4.65%  "
$ws.Range("E33").Value = "  -9.74%  "
$ws.Range("E34").Value = "  -8.76%  "
$ws.Range("D35").Value = "152.08"
$ws.Range("E35").Value = "  -9.93%  "
$ws.Range("D36").Value = "2.67"
$ws.Range("E36").Value = "  -7.49%  "
$ws.Range("D37").Value = "3.13"
$ws.Range("E37").Value = "  +2.15%  "
$ws.Range("E38").Value = "  -7.55%  "
$ws.Range("E39").Value = "  -1.24%  "
$ws.Range("E40").Value = "  -6.56%  "
$ws.Range("E41").Value = "  -11.23%  "
$ws.Range("D42").Value = "3.66"
$ws.Range("E42").Value = "  -9.04%  "
$ws.Range("D43").Value = "0.0320"
$ws.Range("E43").Value = "  -9.57%  "
$ws.Range("D44").Value = "13.77"
$ws.Range("E44").Value = "  +6.22%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "1.725.94"
$ws.Range("E46").Value = "  -7.60%  "
$ws.Range("D48").Value = "84.36"
$ws.Range("E48").Value = "  -14.96%  "
$ws.Range("E49").Value = "  -11.81%  "
$ws.Range("D50").Value = "8.73"
$ws.Range("E50").Value = "  -4.90%  "
$ws.Range("D51").Value = "74.12"
$ws.Range("E51").Value = "  -12.51%  "
